$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new rows of data at the top (rows 1 and 2)
$ws.Range("A1").Value = 1
$ws.Range("B1").Value = 1
$ws.Range("C1").Value = 1

$ws.Range("A2").Value = 2
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 2

# Update the active selection to C3
$ws.Range("C3").Select()
